$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: Find Common Characters (Easy / Arrays-Hashing / Easy) - not yet fully completed
$ws.Cells.Item(46, 1).Value = "Find Common Characters"
$ws.Cells.Item(46, 2).Value = "Easy"
$ws.Cells.Item(46, 3).Value = "Arrays/Hashing"
$ws.Cells.Item(46, 4).Value = "Easy"
$ws.Cells.Item(46, 5).Value = "This one wont be marked as completed until we can do it 100% by ourselves tomorrow"

# Row 47: Intersection of 2 Arrays II (Easy / Arrays-Hashing / Medium) - completed 8/28/2024, no help needed
$ws.Cells.Item(47, 1).Value = "Intersection of 2 Arrays II"
$ws.Cells.Item(47, 2).Value = "Easy"
$ws.Cells.Item(47, 3).Value = "Arrays/Hashing"
$ws.Cells.Item(47, 4).Value = "Medium"
$ws.Cells.Item(47, 5).Value = "This one was done well, not sure why we cant shorten the algorithm though, to just check if nums2[j] is in the dict and is valid"
$ws.Cells.Item(47, 6).Value = 45532
$ws.Cells.Item(47, 6).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(47, 7).Value = "No"
$ws.Cells.Item(47, 9).Value = "Initialize count vector for one str, cmp it to second string, add results if there is connection in between, could also make 2 count vectors and just AND them together, considering that diff values could be possible"

# Recalculate so the COUNTA summary formula in A2 reflects the two newly added rows
$wb.Application.Calculate()

# Update the view/selection to match where the editor ended up (row 22 visible, F45 selected)
$ws.Range("A22").Select() | Out-Null
$ws.Range("F45").Select() | Out-Null
